$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 7412.066666666667, 8011, 6888, 0.1150546471277873),
    @(1, 7914.2, 8561, 7472, 0.1203464428583781),
    @(2, 8134.533333333334, 8709, 7353, 0.1216469605763753),
    @(3, 8661.866666666667, 9485, 7842, 0.1183188676834106),
    @(4, 7755.366666666667, 8157, 7334, 0.1147886435190837),
    @(5, 7961.333333333333, 8621, 7428, 0.1210931936899821),
    @(6, 7914.7, 8429, 7129, 0.1214846213658651),
    @(7, 7682.966666666666, 8354, 6682, 0.1205755710601807),
    @(8, 7146.166666666667, 7599, 6387, 0.1155076026916504),
    @(9, 7523.2, 7979, 6389, 0.1155570189158122)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
